$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "41.475.76"
$ws.Cells.Item(2, 5).Value = "  +4.34%  "
$ws.Cells.Item(3, 4).Value = "2.219.86"
$ws.Cells.Item(3, 5).Value = "  +2.70%  "
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  +0.12%  "
$ws.Cells.Item(5, 4).Value = "'231.31"
$ws.Cells.Item(5, 5).Value = "  +1.99%  "
$ws.Cells.Item(6, 4).Value = "'0.625"
$ws.Cells.Item(6, 5).Value = "  +0.41%  "
$ws.Cells.Item(7, 4).Value = "'61.15"
$ws.Cells.Item(7, 5).Value = "  -2.80%  "
$ws.Cells.Item(8, 5).Value = "  +0.10%  "
$ws.Cells.Item(9, 5).Value = "  +2.83%  "
$ws.Cells.Item(10, 4).Value = "'58.77"
$ws.Cells.Item(10, 5).Value = "  +0.01%  "
$ws.Cells.Item(11, 4).Value = "'0.0889"
$ws.Cells.Item(11, 5).Value = "  +5.74%  "
$ws.Cells.Item(12, 5).Value = "  +0.40%  "
$ws.Cells.Item(13, 4).Value = "2.550.66"
$ws.Cells.Item(13, 5).Value = "  +2.70%  "
$ws.Cells.Item(14, 4).Value = "'15.66"
$ws.Cells.Item(14, 5).Value = "  -1.36%  "
$ws.Cells.Item(15, 4).Value = "'21.76"
$ws.Cells.Item(15, 5).Value = "  -0.14%  "
$ws.Cells.Item(16, 4).Value = "'0.798"
$ws.Cells.Item(16, 5).Value = "  -0.95%  "
$ws.Cells.Item(17, 5).Value = "  +0.67%  "
$ws.Cells.Item(18, 4).Value = "2.213.04"
$ws.Cells.Item(18, 5).Value = "  +2.66%  "
$ws.Cells.Item(19, 4).Value = "41.354.47"
$ws.Cells.Item(19, 5).Value = "  +4.05%  "
$ws.Cells.Item(20, 4).Value = "'72.92"
$ws.Cells.Item(20, 5).Value = "  +1.59%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0896"
$ws.Cells.Item(21, 5).Value = "  +5.89%  "
$ws.Cells.Item(22, 4).Value = "'6.05"
$ws.Cells.Item(22, 5).Value = "  +0.60%  "
$ws.Cells.Item(23, 4).Value = "'250.25"
$ws.Cells.Item(23, 5).Value = "  +9.67%  "
$ws.Cells.Item(24, 5).Value = "  -0.02%  "
$ws.Cells.Item(25, 5).Value = "  +1.16%  "
$ws.Cells.Item(26, 5).Value = "  -0.27%  "
$ws.Cells.Item(27, 4).Value = "'9.44"
$ws.Cells.Item(27, 5).Value = "  +0.09%  "
$ws.Cells.Item(28, 4).Value = "'168.37"
$ws.Cells.Item(28, 5).Value = "  -2.20%  "
$ws.Cells.Item(29, 5).Value = "  -0.46%  "
$ws.Cells.Item(30, 4).Value = "'19.94"
$ws.Cells.Item(30, 5).Value = "  +1.76%  "
$ws.Cells.Item(31, 4).Value = "'1.43"
$ws.Cells.Item(31, 5).Value = "  -1.14%  "
$ws.Cells.Item(32, 4).Value = "'2.64"
$ws.Cells.Item(32, 5).Value = "  -2.11%  "
$ws.Cells.Item(33, 5).Value = "  +0.79%  "
$ws.Cells.Item(34, 4).Value = "'4.98"
$ws.Cells.Item(34, 5).Value = "  +6.24%  "
$ws.Cells.Item(35, 4).Value = "'4.63"
$ws.Cells.Item(35, 5).Value = "  +0.83%  "
$ws.Cells.Item(36, 4).Value = "'0.0626"
$ws.Cells.Item(36, 5).Value = "  +1.50%  "
$ws.Cells.Item(37, 4).Value = "'6.56"
$ws.Cells.Item(37, 5).Value = "  -5.36%  "
$ws.Cells.Item(38, 4).Value = "'3.70"
$ws.Cells.Item(38, 5).Value = "  -1.51%  "
$ws.Cells.Item(40, 4).Value = "'0.000246"
$ws.Cells.Item(40, 5).Value = "  +29.50%  "
$ws.Cells.Item(41, 4).Value = "'1.00"
$ws.Cells.Item(41, 5).Value = "  +0.21%  "
$ws.Cells.Item(42, 2).Value = "FTXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(42, 4).Value = "'4.83"
$ws.Cells.Item(42, 5).Value = "  -1.61%  "
$ws.Cells.Item(43, 2).Value = "VeChain"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(43, 4).Value = "'0.0237"
$ws.Cells.Item(43, 5).Value = "  +4.48%  "
$ws.Cells.Item(44, 4).Value = "'8.60"
$ws.Cells.Item(44, 5).Value = "  +8.44%  "
$ws.Cells.Item(45, 4).Value = "'0.0982"
$ws.Cells.Item(45, 5).Value = "  +6.22%  "
$ws.Cells.Item(46, 4).Value = "'99.04"
$ws.Cells.Item(46, 5).Value = "  -3.59%  "
$ws.Cells.Item(47, 5).Value = "  -0.26%  "
$ws.Cells.Item(48, 4).Value = "1.464.49"
$ws.Cells.Item(48, 5).Value = "  -3.17%  "
$ws.Cells.Item(49, 4).Value = "'16.59"
$ws.Cells.Item(49, 5).Value = "  -5.49%  "
$ws.Cells.Item(50, 5).Value = "  -0.94%  "
$ws.Cells.Item(51, 2).Value = "ARBITRUM"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(51, 4).Value = "'1.08"
$ws.Cells.Item(51, 5).Value = "  -1.18%  "
